$wb = $excel.ActiveWorkbook

# Locate the sheet to use as a template for the header style, and the sheet
# after which the new week's sheet should be inserted.
$templateSheet = $wb.Worksheets.Item("2026-01-19")

# Add the new worksheet right after the most recent existing weekly sheet.
$ws = $wb.Worksheets.Add($null, $templateSheet)
$ws.Name = "2026-01-26"

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "author"
$ws.Range("D1").Value = "latest_episode"

# Copy the header formatting (bold, centered, bordered) from the previous week's
# sheet so the new header matches the established style exactly.
$templateSheet.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Ranking data rows
$data = @(
  @(1, "宇崎ちゃんは遊びたい！", "丈(著者)", "第130話"),
  @(2, "転生したらスライムだった件", "原作：伏瀬 漫画：川上泰樹 キャラクター原案：みっつばー", "★描き下ろし色紙プレゼントのお知らせ★"),
  @(3, "ワンパンマン", "原作/ＯＮＥ 作画/村田雄介", "219撃目"),
  @(4, "異世界おじさん", "殆ど死んでいる(著者)", "番外編10"),
  @(5, "新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる", "漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea", "第75話"),
  @(6, "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第35話：発動の線引き③"),
  @(7, "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第７３話『上昇停止』③"),
  @(8, "勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが", "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)", "第7話 前編"),
  @(9, "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第７１話　いきなり告白の戦いが始まった（２）"),
  @(10, "元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～", "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)", "第82話その2"),
  @(11, "怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった", "菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)", "第18話後半"),
  @(12, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５６話　勇者を蹂躙する器用貧乏（３）"),
  @(13, "実は俺、最強でした？", "原作：澄守 彩 漫画：高橋 愛", "第137話　まずは一戦！・後編"),
  @(14, "え、社内システム全てワンオペしている私を解雇ですか？", "漫画：伊於 原作：下城米雪 キャラクター原案：icchi", "第35話"),
  @(15, "異世界魔王と召喚少女の奴隷魔術", "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大", "第132話　王都に戻ってみるⅠ（後編）"),
  @(16, "まんきつしたい常連さん", "しんみりん(著者)", "第52話後編"),
  @(17, "異世界のんびり農家", "剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)", "第315話"),
  @(18, "バキ外伝 烈海王は異世界転生しても一向にかまわんッッ", "板垣恵介 猪原賽 陸井栄史", "第88話　種"),
  @(19, "クセ強彼女は床にいざなう", "須河篤志(著者)", "休載イラスト"),
  @(20, "ぽんドロイド！ はまさん", "はれやまはれぞう(著者)", "第15話"),
  @(21, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "第２６食　夜の森の焼きリンゴ、パクパクですわ！（３）"),
  @(22, "おとなりのダウナーさんは無理させない", "瑠璃いろ(著者)", "第14話"),
  @(23, "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～", "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき", "第36話 独身貴族は生ハムの原木を買う（2）"),
  @(24, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第77話③"),
  @(25, "「ククク……。奴は四天王の中でも最弱」と解雇された俺、なぜか勇者と聖女の師匠になる", "漫画：芳橋アツシ 原作：延野正行 キャラクター原案：坂野杏梨", "第45話　奴は巡る（前編）"),
  @(26, "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜", "戸賀 環 坂木持丸 riritto", "第60話①　ニコニ坑道の異変を調べてみた"),
  @(27, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第146話 よくわからなけれどシビれるみたいです（２）"),
  @(28, "濁る瞳で何を願う ハイセルク戦記", "トルトネン 創-taro 斎藤八呑", "第37話　意思を持つ天災②"),
  @(29, "アイドル辞めるけど結婚してくれますか!?", "三吉汐美(著者)", "第20話前半"),
  @(30, "理想のヒモ生活", "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)", "第91話　その4"),
  @(31, "バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～", "板垣恵介 林たかあき", "第62話 勝利への執着"),
  @(32, "姫様“拷問”の時間です", "原作:春原ロビンソン　漫画:ひらけい", "拷問167"),
  @(33, "Lv２からチートだった元勇者候補のまったり異世界ライフ", "糸町秋音（漫画） 鬼ノ城ミヤ（原作） 片桐（キャラクター原案）", "第65話　はじめての"),
  @(34, "ふかふかダンジョン攻略記～俺の異世界転生冒険譚～", "KAKERU", "第72話「弓王」（後半）"),
  @(35, "世界最強の後衛 ～迷宮国の新人探索者～", "力蔵(著者) とーわ(原作) 風花風花(キャラクター原案)", "第41話-3"),
  @(36, "バーサス", "原作：ONE 漫画：あずま京太郎 構成：bose", "第30話　勇者の意思（2）"),
  @(37, "老後に備えて異世界で８万枚の金貨を貯めます", "FUNA 東西 モトエ恵介", "第126話　襲撃［その５］"),
  @(38, "回復術士のやり直し", "月夜涙(原作) 羽賀ソウケン(漫画) しおこんぶ(キャラクター原案)", "第75話-3"),
  @(39, "最果てのパラディン", "奥橋睦（漫画） 柳野かなた（原作） 輪くすさが（キャラクター原案）", "第69話　無敵の巨人Ⅱ"),
  @(40, "婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版", "漫画/すたひろ 原作/Y.A", "chapter78【41話①】"),
  @(41, "10年ぶりに再会したクソガキは清純美少女JKに成長していた", "緑青黒羽（漫画） 館西夕木（原作） ひげ猫（キャラクター原案）", "第9話　有月勇のトラウマ（後編）"),
  @(42, "ゲーム悪役貴族に転生した俺は、チート筋肉で無双する", "昼行燈（原作） しいたけ元帥（漫画）", "第37話"),
  @(43, "サーシャちゃんとクラスメイトオタクくん", "はぐはぐ(著者)", "第92話"),
  @(44, "ライブダンジョン！", "ことりりょう(作画) dy冷凍(原作) Mika Pikazo(キャラクター原案)", "第93話前半"),
  @(45, "魔王になったので、ダンジョン造って人外娘とほのぼのする", "遠野ノオト(作画) 流優(原作) だぶ竜(キャラクター原案)", "第86話前半"),
  @(46, "解雇された暗黒兵士(30代)のスローなセカンドライフ", "岡沢六十四 るれくちぇ sage・ジョー", "第77話(前編) レッツ闇鍋パーティー！"),
  @(47, "黄金の経験値", "原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)", "第22話（前編）"),
  @(48, "女子高生の無駄づかい", "ビーノ(著者)", "第140話　かいかく"),
  @(49, "最弱貴族に転生したので悪役たちを集めてみた", "空野進 sorani ファルまろ", "第16話 最弱貴族、魔王を拾う(３)"),
  @(50, "俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。", "狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)", "最終話-1")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

